$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows of hours data
# Shared-string table indices are assigned in the order the distinct text
# values are first written, so set the cell text in that exact order:
#   7 -> "Setting up the ui dynamically..."
#   8 -> "Setting up inventory..."
#   9 -> "Using ui to influence gameplay..."
$ws.Range("A8").Value = "Setting up the ui dynamically and making it organized"
$ws.Range("A10").Value = "Setting up inventory, easilly adding/removing item's and displaying them correctly"
$ws.Range("A9").Value = "Using ui to influence gameplay: equip item's switch character"

$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 3
$ws.Range("B10").Value = 1

# Adjust column width for column A to fit new, longer text.
# (Excel quantizes ColumnWidth to the nearest renderable pixel/character
# increment when it writes the file, so 61.3 is the input that lands on
# the stored width closest to the target of 62.125 "characters".)
$ws.Columns.Item(1).ColumnWidth = 61.3

# Update the selected cell to match the final state
$ws.Range("B8").Select()
